# Updated cryptos list on Mon Apr  3 14:48:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string value to a cell without Excel coercing
# numeric-looking text (e.g. "312.02", "1.002") into a floating point
# number (which would lose the exact printed form / add rounding noise).
# Forcing the NumberFormat to Text ("@") for the write keeps the value an
# exact string, and resetting the format + style afterwards avoids leaving
# a stray custom cell style behind.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.180.01"
Set-TextValue "E2" "  -0.76%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.806.37"
Set-TextValue "E3" "  -0.83%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.18%  "

# Row 5 - BNB
Set-TextValue "D5" "312.02"
Set-TextValue "E5" "  -1.13%  "

# Row 6 - USDC
Set-TextValue "E6" "  +0.00%  "

# Row 7 - XRP
Set-TextValue "D7" "0.5118"
Set-TextValue "E7" "  -2.11%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3921"
Set-TextValue "E8" "  +1.77%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07802"
Set-TextValue "E9" "  -2.99%  "

# Row 10 - Polygon
Set-TextValue "D10" "1.106"
Set-TextValue "E10" "  -0.66%  "

# Row 11 - OKB
Set-TextValue "D11" "41.05"
Set-TextValue "E11" "  -1.96%  "

# Row 12 - Polkadot
Set-TextValue "D12" "6.362"
Set-TextValue "E12" "  -0.47%  "

# Row 13 - BinanceUSD
Set-TextValue "D13" "1.002"
Set-TextValue "E13" "  -0.16%  "

# Row 14 - Solana
Set-TextValue "D14" "20.37"
Set-TextValue "E14" "  -2.58%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.315"
Set-TextValue "E15" "  -1.53%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.804.98"
Set-TextValue "E16" "  -0.67%  "

# Row 17 - Litecoin
Set-TextValue "D17" "92.45"
Set-TextValue "E17" "  -2.09%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.00001075"
Set-TextValue "E18" "  -2.85%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06594"
Set-TextValue "E19" "  -0.71%  "

# Row 20 - Dai
Set-TextValue "D20" "1.003"
Set-TextValue "E20" "  -0.02%  "

# Row 21 - Avalanche
Set-TextValue "D21" "17.28"
Set-TextValue "E21" "  -1.89%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.999"
Set-TextValue "E22" "  -0.30%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "28.204.00"
Set-TextValue "E23" "  -0.86%  "

# Row 24 - Cosmos
Set-TextValue "D24" "11.12"
Set-TextValue "E24" "  -2.10%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.236"
Set-TextValue "E25" "  -0.60%  "

# Row 26 - Monero
Set-TextValue "D26" "160.73"
Set-TextValue "E26" "  +1.11%  "

# Row 27 - LidoDAOToken
Set-TextValue "D27" "2.460"
Set-TextValue "E27" "  +2.09%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextValue "D28" "2.013.26"
Set-TextValue "E28" "  -0.81%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "20.46"
Set-TextValue "E29" "  -1.87%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "127.53"
Set-TextValue "E30" "  +2.48%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.1093"
Set-TextValue "E31" "  -1.52%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "1.057"
Set-TextValue "E32" "  -1.89%  "

# Row 33 - HuobiToken
Set-TextValue "D33" "3.657"
Set-TextValue "E33" "  -0.62%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.563"
Set-TextValue "E34" "  -1.95%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.07075"
Set-TextValue "E35" "  -3.05%  "

# Row 36 - FraxShare
Set-TextValue "D36" "9.138"
Set-TextValue "E36" "  +3.96%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.02345"
Set-TextValue "E37" "  +0.01%  "

# Row 38 - Algorand
Set-TextValue "D38" "0.2171"
Set-TextValue "E38" "  -1.48%  "

# Row 39 - Aptos
Set-TextValue "D39" "11.59"
Set-TextValue "E39" "  -5.24%  "

# Row 40 - InternetComputer(DFINITY)
Set-TextValue "D40" "5.005"
Set-TextValue "E40" "  -1.98%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.6163"
Set-TextValue "E41" "  -2.35%  "

# Row 42 - Frax
Set-TextValue "D42" "1.002"
Set-TextValue "E42" "  +0.11%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "1.158"
Set-TextValue "E43" "  -2.09%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "13.13"
Set-TextValue "E44" "  -2.79%  "

# Rows 45 & 46 - Decentraland and WEMIXTOKEN swap ranking positions
Set-TextValue "B45" "WEMIXTOKEN"
Set-TextValue "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D45" "1.305"
Set-TextValue "E45" "  -5.68%  "

Set-TextValue "B46" "Decentraland"
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5940"
Set-TextValue "E46" "  -3.31%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "3.725"
Set-TextValue "E47" "  -2.07%  "

# Row 48 - Quant
Set-TextValue "D48" "125.20"
Set-TextValue "E48" "  -1.44%  "

# Row 49 - EOS
Set-TextValue "E49" "  -0.91%  "

# Row 50 - NEARProtocol
Set-TextValue "D50" "1.920"
Set-TextValue "E50" "  -2.46%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.06775"
Set-TextValue "E51" "  -1.85%  "
